$d = $word.ActiveDocument
$d.Content.Font.Name = "Arial"
